$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column C ("תיאור"), shifting
# the existing C:F columns (תיאור, מספר חשבונית, פטור ממעמ, סטטוס) to D:G.
$ws.Columns("C:C").Insert()

# New column header + the renamed "מספר חשבונית" -> "מספר מסמך" header
# (now in column E after the shift).
$ws.Range("C1").Value = "סוג מסמך"
$ws.Range("E1").Value = "מספר מסמך"

# New column values: document type per expense row.
$ws.Range("C2").Value = "tax_invoice"
$ws.Range("C3").Value = "tax_invoice_receipt"
$ws.Range("C4").Value = "tax_invoice"
$ws.Range("C5").Value = "receipt"

# Row 3 ("שכירות משרד") status flips from pending to paid.
$ws.Range("G3").Value = "paid"
